# Update market-price-derived Leve profit columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 561.6667
$ws.Cells.Item(12, 9).Value = 597.5
$ws.Cells.Item(12, 10).Value = 490
$ws.Cells.Item(12, 11).Value = 597.5
$ws.Cells.Item(12, 12).Value = 490
$ws.Cells.Item(12, 13).Value = -427.5
$ws.Cells.Item(12, 14).Value = -830
# Row 40
$ws.Cells.Item(40, 8).Value = 2558.0186
$ws.Cells.Item(40, 9).Value = 1500
$ws.Cells.Item(40, 10).Value = 8641.625
$ws.Cells.Item(40, 11).Value = 1500
$ws.Cells.Item(40, 12).Value = 8641.625
$ws.Cells.Item(40, 13).Value = -1325
$ws.Cells.Item(40, 14).Value = -8991.625
# Row 62
$ws.Cells.Item(62, 8).Value = 13273.521
$ws.Cells.Item(62, 9).Value = 17613.928
$ws.Cells.Item(62, 10).Value = 6521.778
$ws.Cells.Item(62, 11).Value = 17613.928
$ws.Cells.Item(62, 12).Value = 6521.778
$ws.Cells.Item(62, 13).Value = -16989.928
$ws.Cells.Item(62, 14).Value = -7769.778
# Row 65
$ws.Cells.Item(65, 8).Value = 13273.521
$ws.Cells.Item(65, 9).Value = 17613.928
$ws.Cells.Item(65, 10).Value = 6521.778
$ws.Cells.Item(65, 11).Value = 88069.64
$ws.Cells.Item(65, 12).Value = 32608.89
$ws.Cells.Item(65, 13).Value = -84949.64
$ws.Cells.Item(65, 14).Value = -38848.89
# Row 116
$ws.Cells.Item(116, 8).Value = 4481.8887
$ws.Cells.Item(116, 9).Value = 4788.727
$ws.Cells.Item(116, 11).Value = 4788.727
$ws.Cells.Item(116, 13).Value = -1346.727
# Row 132
$ws.Cells.Item(132, 8).Value = 3043.4473
$ws.Cells.Item(132, 9).Value = 1860.3793
$ws.Cells.Item(132, 10).Value = 6855.5557
$ws.Cells.Item(132, 11).Value = 5581.1379
$ws.Cells.Item(132, 12).Value = 20566.6671
$ws.Cells.Item(132, 13).Value = -3051.1379
$ws.Cells.Item(132, 14).Value = -25626.6671
# Row 136
$ws.Cells.Item(136, 8).Value = 37045.6
$ws.Cells.Item(136, 10).Value = 37045.6
$ws.Cells.Item(136, 12).Value = 37045.6
$ws.Cells.Item(136, 14).Value = -47245.6
# Row 138
$ws.Cells.Item(138, 8).Value = 2487.5852
$ws.Cells.Item(138, 9).Value = 1237.8983
$ws.Cells.Item(138, 10).Value = 4594.2
$ws.Cells.Item(138, 11).Value = 3713.6949
$ws.Cells.Item(138, 12).Value = 13782.6
$ws.Cells.Item(138, 13).Value = 1426.3051
$ws.Cells.Item(138, 14).Value = -24062.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 960.61536
$ws.Cells.Item(45, 9).Value = 957.5
$ws.Cells.Item(45, 10).Value = 965.6
$ws.Cells.Item(45, 11).Value = 957.5
$ws.Cells.Item(45, 12).Value = 965.6
$ws.Cells.Item(45, 13).Value = -580.5
$ws.Cells.Item(45, 14).Value = -1719.6
# Row 74
$ws.Cells.Item(74, 8).Value = 187856.39
$ws.Cells.Item(74, 9).Value = 213601.4
$ws.Cells.Item(74, 10).Value = 87021.75
$ws.Cells.Item(74, 11).Value = 213601.4
$ws.Cells.Item(74, 12).Value = 87021.75
$ws.Cells.Item(74, 13).Value = -212727.4
$ws.Cells.Item(74, 14).Value = -88769.75
# Row 77
$ws.Cells.Item(77, 8).Value = 187856.39
$ws.Cells.Item(77, 9).Value = 213601.4
$ws.Cells.Item(77, 10).Value = 87021.75
$ws.Cells.Item(77, 11).Value = 1068007
$ws.Cells.Item(77, 12).Value = 435108.75
$ws.Cells.Item(77, 13).Value = -1063639
$ws.Cells.Item(77, 14).Value = -443844.75
# Row 102
$ws.Cells.Item(102, 8).Value = 11246.667
$ws.Cells.Item(102, 9).Value = 1500
$ws.Cells.Item(102, 10).Value = 20993.334
$ws.Cells.Item(102, 11).Value = 1500
$ws.Cells.Item(102, 12).Value = 20993.334
$ws.Cells.Item(102, 13).Value = 122
$ws.Cells.Item(102, 14).Value = -24237.334
# Row 132
$ws.Cells.Item(132, 8).Value = 19636.967
$ws.Cells.Item(132, 9).Value = 26643.537
$ws.Cells.Item(132, 11).Value = 79930.611
$ws.Cells.Item(132, 13).Value = -77400.611

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2991.6875
$ws.Cells.Item(31, 9).Value = 2131.743
$ws.Cells.Item(31, 10).Value = 5306.923
$ws.Cells.Item(31, 11).Value = 2131.743
$ws.Cells.Item(31, 12).Value = 5306.923
$ws.Cells.Item(31, 13).Value = -1836.743
$ws.Cells.Item(31, 14).Value = -5896.923
# Row 34
$ws.Cells.Item(34, 8).Value = 2991.6875
$ws.Cells.Item(34, 9).Value = 2131.743
$ws.Cells.Item(34, 10).Value = 5306.923
$ws.Cells.Item(34, 11).Value = 2131.743
$ws.Cells.Item(34, 12).Value = 5306.923
$ws.Cells.Item(34, 13).Value = -1929.743
$ws.Cells.Item(34, 14).Value = -5710.923
# Row 134
$ws.Cells.Item(134, 8).Value = 1306.1333
$ws.Cells.Item(134, 9).Value = 888.9796
$ws.Cells.Item(134, 10).Value = 3164.3635
$ws.Cells.Item(134, 11).Value = 2666.9388
$ws.Cells.Item(134, 12).Value = 9493.0905
$ws.Cells.Item(134, 13).Value = -131.9387999999999
$ws.Cells.Item(134, 14).Value = -14563.0905

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Cells.Item(56, 8).Value = 5453.636
$ws.Cells.Item(56, 9).Value = 5453.636
$ws.Cells.Item(56, 11).Value = 5453.636
$ws.Cells.Item(56, 13).Value = -4923.636
# Row 64
$ws.Cells.Item(64, 8).Value = 2370.3333
$ws.Cells.Item(64, 9).Value = 1409.1666
$ws.Cells.Item(64, 10).Value = 3011.111
$ws.Cells.Item(64, 11).Value = 4227.4998
$ws.Cells.Item(64, 12).Value = 9033.332999999999
$ws.Cells.Item(64, 13).Value = -3957.4998
$ws.Cells.Item(64, 14).Value = -9573.332999999999
# Row 67
$ws.Cells.Item(67, 8).Value = 2370.3333
$ws.Cells.Item(67, 9).Value = 1409.1666
$ws.Cells.Item(67, 10).Value = 3011.111
$ws.Cells.Item(67, 11).Value = 4227.4998
$ws.Cells.Item(67, 12).Value = 9033.332999999999
$ws.Cells.Item(67, 13).Value = -3291.4998
$ws.Cells.Item(67, 14).Value = -10905.333
# Row 70
$ws.Cells.Item(70, 8).Value = 3929.4443
$ws.Cells.Item(70, 9).Value = 475
$ws.Cells.Item(70, 10).Value = 4361.25
$ws.Cells.Item(70, 11).Value = 1425
$ws.Cells.Item(70, 12).Value = 13083.75
$ws.Cells.Item(70, 13).Value = -1110
$ws.Cells.Item(70, 14).Value = -13713.75
# Row 73
$ws.Cells.Item(73, 8).Value = 3929.4443
$ws.Cells.Item(73, 9).Value = 475
$ws.Cells.Item(73, 10).Value = 4361.25
$ws.Cells.Item(73, 11).Value = 1425
$ws.Cells.Item(73, 12).Value = 13083.75
$ws.Cells.Item(73, 13).Value = -333
$ws.Cells.Item(73, 14).Value = -15267.75
# Row 87
$ws.Cells.Item(87, 8).Value = 5162.8
$ws.Cells.Item(87, 9).Value = 3638
$ws.Cells.Item(87, 10).Value = 7450
$ws.Cells.Item(87, 11).Value = 10914
$ws.Cells.Item(87, 12).Value = 22350
$ws.Cells.Item(87, 13).Value = -9666
$ws.Cells.Item(87, 14).Value = -24846
# Row 90
$ws.Cells.Item(90, 8).Value = 5162.8
$ws.Cells.Item(90, 9).Value = 3638
$ws.Cells.Item(90, 10).Value = 7450
$ws.Cells.Item(90, 11).Value = 32742
$ws.Cells.Item(90, 12).Value = 67050
$ws.Cells.Item(90, 13).Value = -26502
$ws.Cells.Item(90, 14).Value = -79530
# Row 98
$ws.Cells.Item(98, 8).Value = 6231.091
$ws.Cells.Item(98, 10).Value = 8430.25
$ws.Cells.Item(98, 12).Value = 25290.75
$ws.Cells.Item(98, 14).Value = -28286.75
# Row 114
$ws.Cells.Item(114, 8).Value = 353.4
$ws.Cells.Item(114, 9).Value = 291.75
$ws.Cells.Item(114, 10).Value = 600
$ws.Cells.Item(114, 11).Value = 875.25
$ws.Cells.Item(114, 12).Value = 1800
$ws.Cells.Item(114, 13).Value = 2378.75
$ws.Cells.Item(114, 14).Value = -8308
# Row 131
$ws.Cells.Item(131, 8).Value = 1484.7455
$ws.Cells.Item(131, 10).Value = 1425.6459
$ws.Cells.Item(131, 12).Value = 4276.9377
$ws.Cells.Item(131, 14).Value = -14356.9377

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 5461.579
$ws.Cells.Item(107, 9).Value = 9195
$ws.Cells.Item(107, 11).Value = 9195
$ws.Cells.Item(107, 13).Value = -7275
# Row 122
$ws.Cells.Item(122, 8).Value = 251704.67
$ws.Cells.Item(122, 9).Value = 274241.47
$ws.Cells.Item(122, 10).Value = 3800
$ws.Cells.Item(122, 11).Value = 822724.4099999999
$ws.Cells.Item(122, 12).Value = 11400
$ws.Cells.Item(122, 13).Value = -820274.4099999999
$ws.Cells.Item(122, 14).Value = -16300
# Row 132
$ws.Cells.Item(132, 8).Value = 3174.0386
$ws.Cells.Item(132, 9).Value = 2964.7026
$ws.Cells.Item(132, 10).Value = 3690.4
$ws.Cells.Item(132, 11).Value = 8894.1078
$ws.Cells.Item(132, 12).Value = 11071.2
$ws.Cells.Item(132, 13).Value = -6364.1078
$ws.Cells.Item(132, 14).Value = -16131.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1717.3334
$ws.Cells.Item(16, 9).Value = 1596.8462
$ws.Cells.Item(16, 10).Value = 2164.8572
$ws.Cells.Item(16, 11).Value = 1596.8462
$ws.Cells.Item(16, 12).Value = 2164.8572
$ws.Cells.Item(16, 13).Value = -1426.8462
$ws.Cells.Item(16, 14).Value = -2504.8572
# Row 47
$ws.Cells.Item(47, 8).Value = 17857.5
$ws.Cells.Item(47, 10).Value = 17857.5
$ws.Cells.Item(47, 12).Value = 17857.5
$ws.Cells.Item(47, 14).Value = -18837.5
# Row 52
$ws.Cells.Item(52, 8).Value = 17857.5
$ws.Cells.Item(52, 10).Value = 17857.5
$ws.Cells.Item(52, 12).Value = 17857.5
$ws.Cells.Item(52, 14).Value = -18323.5
# Row 61
$ws.Cells.Item(61, 8).Value = 1038.5
$ws.Cells.Item(61, 9).Value = 606.8889
$ws.Cells.Item(61, 10).Value = 2333.3333
$ws.Cells.Item(61, 11).Value = 606.8889
$ws.Cells.Item(61, 12).Value = 2333.3333
$ws.Cells.Item(61, 13).Value = -404.8889
$ws.Cells.Item(61, 14).Value = -2737.3333
# Row 113
$ws.Cells.Item(113, 8).Value = 1038.5
$ws.Cells.Item(113, 9).Value = 606.8889
$ws.Cells.Item(113, 10).Value = 2333.3333
$ws.Cells.Item(113, 11).Value = 606.8889
$ws.Cells.Item(113, 12).Value = 2333.3333
$ws.Cells.Item(113, 13).Value = 1563.1111
$ws.Cells.Item(113, 14).Value = -6673.3333

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 1256.2031
$ws.Cells.Item(132, 9).Value = 711.32074
$ws.Cells.Item(132, 11).Value = 2133.96222
$ws.Cells.Item(132, 13).Value = 396.0377800000001
